$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.372.09"
$ws.Range("D3").Value = "3.506.40"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.00%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.19"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E8").Value = "  +0.51%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.63"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +5.89%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("E11").Value = "  +3.95%  "
$ws.Range("D12").Value = "4.104.67"
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "3.500.16"
$ws.Range("E15").Value = "  +0.40%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.81"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("D17").Value = "64.363.54"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("E20").Value = "  -0.72%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "392.76"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").Value = "3.646.56"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  +0.48%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.17%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.66"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  +3.61%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.02"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +2.11%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.45"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.12%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +2.28%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.21"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +0.00%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.48"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -5.36%  "
$ws.Range("D34").Value = "3.534.84"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("E35").Value = "  +0.01%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.41"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("E39").Value = "  +2.22%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.28"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +2.50%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0789"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  +1.54%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.06"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("E48").Value = "  +1.15%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.919"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").Value = "2.390.92"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("E51").Value = "  +0.52%  "
